# Weekly refresh of Fruit/Vegetable prices ("Coco") data rows.
# The underlying per-record fields (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg) get
# reshuffled across the existing data rows (2-35) while the rest of each
# row (market/product descriptive columns) stays put.
#
# Build the row->row mapping: data that should END UP on a given row comes
# FROM the row indicated below (keys/values are spreadsheet row numbers).
$map = @{2=3; 3=7; 4=29; 5=21; 6=22; 7=9; 8=8; 9=16; 10=33; 11=18; 12=15; 13=5; 14=2; 15=17; 16=32; 17=35; 18=6; 19=11; 20=24; 21=19; 22=31; 23=23; 24=20; 25=10; 26=27; 27=12; 28=28; 29=34; 30=13; 31=14; 32=25; 33=30; 34=26; 35=4}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "record": D, L, M, N, O, P, R, S
$cols = @(4, 12, 13, 14, 15, 16, 18, 19)

# First, snapshot the current (pre-edit) values for every affected column on
# every data row, so that writing new values does not clobber data that is
# still needed as a source for another row later on.
$snapshot = @{}
foreach ($col in $cols) {
    $colData = @{}
    for ($r = 2; $r -le 35; $r++) {
        $colData[$r] = $ws.Cells.Item($r, $col).Value2
    }
    $snapshot[$col] = $colData
}

# Now write the shuffled values back out according to the mapping.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($col in $cols) {
        $colData = $snapshot[$col]
        $ws.Cells.Item($destRow, $col).Value = $colData[$srcRow]
    }
}
